$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 11111111
$ws.Range("B2").Value = 5

# Remove rows 3 and 4 entirely (data no longer needed)
$ws.Range("A3:B4").Delete()
